# Update countries & provincias Spain
# Refreshes the COVID data snapshot in the "Pais" sheet:
#   - Kazajistan overtakes Emiratos Arabes Unidos (rows 34/35 swap + new Kazajistan data)
#   - Fiyi's case count moves it ahead of Curazao (rows 199-203 cascade)
#   - Honduras (row 55) and Birmania (row 164) values refreshed in place
#   - "Datos actualizados" timestamp bumped from 04:32 to 05:49

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 05:49"

# --- Rows 34/35: Kazajistan now has more cases than Emiratos Arabes Unidos ---
$ws.Range("A34").Value = "Kazajistan"
$ws.Range("B34").Value = 54747
$ws.Range("C34").Value = 1726
$ws.Range("D34").Value = 35137
$ws.Range("E34").Value = 19346
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 264

$ws.Range("A35").Value = "Emiratos Arabes Unidos"
$ws.Range("B35").Value = 53577
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 43570
$ws.Range("E35").Value = 9679
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 328

# --- Row 55: Honduras refresh ---
$ws.Range("B55").Value = 26384
$ws.Range("C55").Value = 406
$ws.Range("D55").Value = 2779
$ws.Range("E55").Value = 22901
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 10
$ws.Range("H55").Value = 704

# --- Row 164: Birmania refresh ---
$ws.Range("B164").Value = 321
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 254
$ws.Range("E164").Value = 61
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 6

# --- Rows 199-203: Fiyi jumps ahead of Curazao, Timor Oriental, Granada, Santa Lucia ---
$ws.Range("A199").Value = "Fiyi"
$ws.Range("B199").Value = 26
$ws.Range("C199").Value = 5
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 8
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Curazao"
$ws.Range("B200").Value = 25
$ws.Range("C200").Value = 2
$ws.Range("D200").Value = 24
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1

$ws.Range("A201").Value = "Timor Oriental"
$ws.Range("B201").Value = 24
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 24
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

$ws.Range("A202").Value = "Granada"
$ws.Range("B202").Value = 23
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 23
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 22
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 19
$ws.Range("E203").Value = 3
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0
